# giants-gold (Version 2): move the "Meta description" paragraph from the
# top of the document down to just above the final "Prompt:" paragraph,
# retitling it to the page's H1 text, and turn the old "Prompt:" paragraph
# into the new meta-description blurb.
#
# Net effect (matches the supplied unified diff):
#   1. The paragraph
#        [empty run] + "Meta description" (bold) + ": Find out why ... free!"
#      disappears from its original spot (right after the H1 heading).
#   2. A paragraph
#        [empty run] + "Play Giant's Gold Free: A Review of the Unique
#        Two-Grid Game" (bold)
#      appears immediately before the last paragraph of the document (the
#      one that used to start with "Prompt: Can you create a feature
#      image...").
#   3. That last paragraph's text is replaced with "Find out why Giant's
#      Gold is a refreshing change from regular slot games with our
#      review. Play now for free!" while keeping its italic run formatting.

$d = $word.ActiveDocument

# --- Steps 1 & 2: relocate the "Meta description" paragraph -----------
# Cut it from its original location (2nd paragraph of the document) and
# paste it back in immediately before the final paragraph. This is a
# single logical "move", so it naturally produces both the deletion near
# the top and the insertion near the bottom described by the diff.
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Cut()

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertionPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$insertionPoint.Paste()

# The pasted paragraph is now the second-to-last paragraph in the document.
$movedPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)

# Drop the old ": Find out why ... free!" run - only the bold label run
# should remain in this paragraph.
$suffixRange = $d.Range($movedPara.Range.Start, $movedPara.Range.End)
$suffixFound = $suffixRange.Find.Execute(": Find out why Giant's Gold is a refreshing change from regular slot games with our review. Play now for free!")
if ($suffixFound) {
    $suffixRange.Text = ""
}

# Retitle the remaining bold run from "Meta description" to the new
# heading text (kept bold, matching the diff).
$titleRange = $d.Range($movedPara.Range.Start, $movedPara.Range.End)
$titleFound = $titleRange.Find.Execute("Meta description")
if ($titleFound) {
    $titleRange.Text = "Play Giant's Gold Free: A Review of the Unique Two-Grid Game"
}

# --- Step 3: rewrite the old "Prompt:" paragraph -----------------------
# Plain Range.Text assignment (rather than Find/Replace) keeps the
# straight apostrophe intact and preserves the existing <w:i/> run
# formatting of that paragraph.
$promptPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$promptRange = $d.Range($promptPara.Range.Start, $promptPara.Range.End)
$promptRange.Text = "Find out why Giant's Gold is a refreshing change from regular slot games with our review. Play now for free!"
